# "bug in indexing the impact assessment"
#
# Before: row 1 = metric headers, row 2 = the single data row
#         (A2/B2 = "Pumpfed Irrigation for Maize" / "baseline", C2:X2 =
#         the 22 metric values).
#
# After:  a units row is inserted right under the header (new row 2:
#         C2:X2 carry the unit for each metric, e.g. "M GHS", "Mm3", ...),
#         row 3 stays blank, and the data moves down to row 4 (A4/B4 keep
#         the labels that used to sit in A2/B2). While re-indexing which
#         metric fed which column, three values per Saving/Investment/
#         Total-Impact block (the Water/Emission/Land columns: G,H,I and
#         M,N,O and S,T,U) turned out to have been picked up wrong, so
#         they get corrected values in the same edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the old row-2 data down to row 4 -----------------------------
# Formats first (so A4/B4 pick up the bordered/centered style A2/B2 had),
# then values explicitly (several of the numbers change along the way).
$ws.Range("A2:B2").Copy()
$ws.Range("A4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A4").Value = "Pumpfed Irrigation for Maize"
$ws.Range("B4").Value = "baseline"

$ws.Range("C4").Value = 125.8672760472691
$ws.Range("D4").Value = 715.2023191869084
$ws.Range("E4").Value = 5.682194305359529
$ws.Range("F4").Value = 0.1759883499683891
$ws.Range("G4").Value = 7.192313248880737
$ws.Range("H4").Value = 33.94489230068575
$ws.Range("I4").Value = 337.0500250131554
$ws.Range("J4").Value = 39.66921853806707
$ws.Range("K4").Value = 377.7461571722815
$ws.Range("L4").Value = 270.477081390447
$ws.Range("M4").Value = 0.002309998559894666
$ws.Range("N4").Value = 15.77667327413656
$ws.Range("O4").Value = 0.003761469974051579
$ws.Range("P4").Value = 102.3345925394533
$ws.Range("Q4").Value = 0.06649889930849895
$ws.Range("R4").Value = 0.3105979636602569
$ws.Range("S4").Value = -71.92082249024747
$ws.Range("T4").Value = -323.672249732721
$ws.Range("U4").Value = -3370.49648866158
$ws.Range("V4").Value = -294.3575928412174
$ws.Range("W4").Value = -2704.704315005161
$ws.Range("X4").Value = -3777.150973759155

# --- Old row 2 becomes the (blank-labelled) units row -----------------
# A2 had nothing above it in column A going forward (the new units row
# has no label there at all) -> drop the cell (value + formatting) so it
# no longer appears in the sheet, same as column A in row 3 which stays
# completely untouched/absent.
$ws.Range("A2").Clear()
# B2 keeps its bordered style, just loses the "baseline" text.
$ws.Range("B2").ClearContents()
$ws.Range("C2:X2").ClearContents()

# B1 gets the same bordered/centered style as the rest of row 1 (it
# stays textless, sitting above the now-blank B2).
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C2").Value = "M GHS"
$ws.Range("D2").Value = "M GHS"
$ws.Range("E2").Value = "1/years"
$ws.Range("F2").Value = "years"
$ws.Range("G2").Value = "Mm3"
$ws.Range("H2").Value = "kton"
$ws.Range("I2").Value = "kha"
$ws.Range("J2").Value = "M GHS"
$ws.Range("K2").Value = "M GHS"
$ws.Range("L2").Value = "M GHS"
$ws.Range("M2").Value = "Mm3"
$ws.Range("N2").Value = "kton"
$ws.Range("O2").Value = "kha"
$ws.Range("P2").Value = "M GHS"
$ws.Range("Q2").Value = "M GHS"
$ws.Range("R2").Value = "M GHS"
$ws.Range("S2").Value = "Mm3"
$ws.Range("T2").Value = "kton"
$ws.Range("U2").Value = "kha"
$ws.Range("V2").Value = "M GHS"
$ws.Range("W2").Value = "M GHS"
$ws.Range("X2").Value = "M GHS"

# Re-apply the bordered/centered style to the unit cells (the Value=
# assignments above only touch content, not formatting).
$ws.Range("C1").Copy()
$ws.Range("B2:X2").PasteSpecial(-4122)   # xlPasteFormats
